# Refresh profit-calc sheets: update Leve price/profit columns (H:N) with
# newly-fetched market-board averages, as produced by the scheduled pricing runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 767.7143
$ws.Range("J17").Value = 776.5185
$ws.Range("L17").Value = 2329.5555
$ws.Range("N17").Value = -2665.5555

$ws.Range("H64").Value = 3102.963
$ws.Range("I64").Value = 2771.4285
$ws.Range("J64").Value = 3219
$ws.Range("K64").Value = 2771.4285
$ws.Range("L64").Value = 3219
$ws.Range("M64").Value = -2523.4285
$ws.Range("N64").Value = -3715

$ws.Range("H67").Value = 3102.963
$ws.Range("I67").Value = 2771.4285
$ws.Range("J67").Value = 3219
$ws.Range("K67").Value = 2771.4285
$ws.Range("L67").Value = 3219
$ws.Range("M67").Value = -1913.4285
$ws.Range("N67").Value = -4935

$ws.Range("H74").Value = 3006.6667
$ws.Range("I74").Value = 2527.7778
$ws.Range("J74").Value = 3725
$ws.Range("K74").Value = 2527.7778
$ws.Range("L74").Value = 3725
$ws.Range("M74").Value = -1591.7778
$ws.Range("N74").Value = -5597

$ws.Range("H77").Value = 3006.6667
$ws.Range("I77").Value = 2527.7778
$ws.Range("J77").Value = 3725
$ws.Range("K77").Value = 12638.889
$ws.Range("L77").Value = 18625
$ws.Range("M77").Value = -7958.888999999999
$ws.Range("N77").Value = -27985

$ws.Range("H116").Value = 2000
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1442
$ws.Range("N116").ClearContents()

$ws.Range("H129").Value = 18581.63
$ws.Range("J129").Value = 25025.738
$ws.Range("L129").Value = 75077.21400000001
$ws.Range("N129").Value = -85077.21400000001

$ws.Range("H137").Value = 1487.8889
$ws.Range("I137").Value = 1221.1765
$ws.Range("J137").Value = 1941.3
$ws.Range("K137").Value = 3663.5295
$ws.Range("L137").Value = 5823.9
$ws.Range("M137").Value = -1113.5295
$ws.Range("N137").Value = -10923.9


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 694.4483
$ws.Range("I74").Value = 701.3929000000001
$ws.Range("J74").Value = 500
$ws.Range("K74").Value = 701.3929000000001
$ws.Range("L74").Value = 500
$ws.Range("M74").Value = 172.6070999999999
$ws.Range("N74").Value = -2248

$ws.Range("H77").Value = 694.4483
$ws.Range("I77").Value = 701.3929000000001
$ws.Range("J77").Value = 500
$ws.Range("K77").Value = 3506.9645
$ws.Range("L77").Value = 2500
$ws.Range("M77").Value = 861.0355
$ws.Range("N77").Value = -11236

$ws.Range("H88").Value = 915429.6
$ws.Range("I88").Value = 1671767.5
$ws.Range("J88").Value = 7824.2
$ws.Range("K88").Value = 1671767.5
$ws.Range("L88").Value = 7824.2
$ws.Range("M88").Value = -1671361.5
$ws.Range("N88").Value = -8636.200000000001

$ws.Range("H91").Value = 915429.6
$ws.Range("I91").Value = 1671767.5
$ws.Range("J91").Value = 7824.2
$ws.Range("K91").Value = 1671767.5
$ws.Range("L91").Value = 7824.2
$ws.Range("M91").Value = -1670363.5
$ws.Range("N91").Value = -10632.2

$ws.Range("H122").Value = 1308.3914
$ws.Range("I122").Value = 1308.3914
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3925.1742
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1475.1742
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 5147.3237
$ws.Range("I132").Value = 6480.857
$ws.Range("K132").Value = 19442.571
$ws.Range("M132").Value = -16912.571


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 358.13635
$ws.Range("I22").Value = 386.25
$ws.Range("J22").Value = 283.16666
$ws.Range("K22").Value = 386.25
$ws.Range("L22").Value = 283.16666
$ws.Range("M22").Value = -213.25
$ws.Range("N22").Value = -629.16666

$ws.Range("H86").Value = 2896.2307
$ws.Range("I86").Value = 2029.5294
$ws.Range("K86").Value = 2029.5294
$ws.Range("M86").Value = -906.5293999999999

$ws.Range("H89").Value = 2896.2307
$ws.Range("I89").Value = 2029.5294
$ws.Range("K89").Value = 10147.647
$ws.Range("M89").Value = -4531.646999999999

$ws.Range("H94").Value = 1882.1111
$ws.Range("I94").Value = 1527.25
$ws.Range("J94").Value = 2166
$ws.Range("K94").Value = 1527.25
$ws.Range("L94").Value = 2166
$ws.Range("M94").Value = -1076.25
$ws.Range("N94").Value = -3068

$ws.Range("H134").Value = 59268.168
$ws.Range("I134").Value = 129340.875
$ws.Range("J134").Value = 3210
$ws.Range("K134").Value = 388022.625
$ws.Range("L134").Value = 9630
$ws.Range("M134").Value = -385487.625
$ws.Range("N134").Value = -14700


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 37039660
$ws.Range("I62").Value = 2498.75
$ws.Range("J62").Value = 66669388
$ws.Range("K62").Value = 2498.75
$ws.Range("L62").Value = 66669388
$ws.Range("M62").Value = -1874.75
$ws.Range("N62").Value = -66670636

$ws.Range("H65").Value = 37039660
$ws.Range("I65").Value = 2498.75
$ws.Range("J65").Value = 66669388
$ws.Range("K65").Value = 12493.75
$ws.Range("L65").Value = 333346940
$ws.Range("M65").Value = -9373.75
$ws.Range("N65").Value = -333353180


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 23817998
$ws.Range("I9").Value = 6250
$ws.Range("J9").Value = 27786624
$ws.Range("K9").Value = 18750
$ws.Range("L9").Value = 83359872
$ws.Range("M9").Value = -18526
$ws.Range("N9").Value = -83360320

$ws.Range("H33").Value = 7505.0835
$ws.Range("I33").Value = 12747.143
$ws.Range("J33").Value = 166.2
$ws.Range("K33").Value = 76482.85800000001
$ws.Range("L33").Value = 997.1999999999999
$ws.Range("M33").Value = -76199.85800000001
$ws.Range("N33").Value = -1563.2

$ws.Range("H41").Value = 950
$ws.Range("J41").Value = 950
$ws.Range("L41").Value = 2850
$ws.Range("N41").Value = -3526

$ws.Range("H69").Value = 500
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 500
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H93").Value = 2375
$ws.Range("J93").Value = 2375
$ws.Range("L93").Value = 7125
$ws.Range("N93").Value = -10869

$ws.Range("H113").Value = 493.5926
$ws.Range("I113").Value = 496.5
$ws.Range("J113").Value = 487.77777
$ws.Range("K113").Value = 1489.5
$ws.Range("L113").Value = 1463.33331
$ws.Range("M113").Value = 680.5
$ws.Range("N113").Value = -5803.33331

$ws.Range("H131").Value = 1854400.8
$ws.Range("I131").Value = 4992.5
$ws.Range("J131").Value = 3087339.5
$ws.Range("K131").Value = 14977.5
$ws.Range("L131").Value = 9262018.5
$ws.Range("M131").Value = -9937.5
$ws.Range("N131").Value = -9272098.5

$ws.Range("H132").Value = 1759.2307
$ws.Range("I132").Value = 1055.8462
$ws.Range("J132").Value = 2462.6155
$ws.Range("K132").Value = 9502.6158
$ws.Range("L132").Value = 22163.5395
$ws.Range("M132").Value = -6972.6158
$ws.Range("N132").Value = -27223.5395


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 739.2
$ws.Range("I97").Value = 652.9231
$ws.Range("J97").Value = 1300
$ws.Range("K97").Value = 652.9231
$ws.Range("L97").Value = 1300
$ws.Range("M97").Value = -156.9231
$ws.Range("N97").Value = -2292

$ws.Range("H113").Value = 83334260
$ws.Range("I113").Value = 125000650
$ws.Range("J113").Value = 1480
$ws.Range("K113").Value = 125000650
$ws.Range("L113").Value = 1480
$ws.Range("M113").Value = -124998480
$ws.Range("N113").Value = -5820

$ws.Range("H122").Value = 20002256
$ws.Range("I122").Value = 66669500
$ws.Range("J122").Value = 2008.2285
$ws.Range("K122").Value = 200008500
$ws.Range("L122").Value = 6024.6855
$ws.Range("M122").Value = -200006050
$ws.Range("N122").Value = -10924.6855

$ws.Range("H132").Value = 3560.111
$ws.Range("I132").Value = 3233.6667
$ws.Range("K132").Value = 9701.000100000001
$ws.Range("M132").Value = -7171.000100000001


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1786.875
$ws.Range("I16").Value = 682.5
$ws.Range("J16").Value = 5100
$ws.Range("K16").Value = 682.5
$ws.Range("L16").Value = 5100
$ws.Range("M16").Value = -512.5
$ws.Range("N16").Value = -5440

$ws.Range("H68").Value = 1294.3
$ws.Range("I68").Value = 1140
$ws.Range("J68").Value = 1397.1666
$ws.Range("K68").Value = 1140
$ws.Range("L68").Value = 1397.1666
$ws.Range("M68").Value = -391
$ws.Range("N68").Value = -2895.1666

$ws.Range("H71").Value = 1294.3
$ws.Range("I71").Value = 1140
$ws.Range("J71").Value = 1397.1666
$ws.Range("K71").Value = 5700
$ws.Range("L71").Value = 6985.833000000001
$ws.Range("M71").Value = -1956
$ws.Range("N71").Value = -14473.833

$ws.Range("H93").Value = 1503068.8
$ws.Range("I93").Value = 2253719.8
$ws.Range("J93").Value = 1766.6666
$ws.Range("K93").Value = 2253719.8
$ws.Range("L93").Value = 1766.6666
$ws.Range("M93").Value = -2252471.8
$ws.Range("N93").Value = -4262.6666

$ws.Range("H122").Value = 3072.5908
$ws.Range("I122").Value = 3357.6667
$ws.Range("J122").Value = 2730.5
$ws.Range("K122").Value = 10073.0001
$ws.Range("L122").Value = 8191.5
$ws.Range("M122").Value = -7623.000100000001
$ws.Range("N122").Value = -13091.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 458.1875
$ws.Range("I107").Value = 430.4
$ws.Range("J107").Value = 470.81818
$ws.Range("K107").Value = 1291.2
$ws.Range("L107").Value = 1412.45454
$ws.Range("M107").Value = 628.8000000000002
$ws.Range("N107").Value = -5252.45454

$ws.Range("H132").Value = 1239.1666
$ws.Range("I132").Value = 1024.4103
$ws.Range("J132").Value = 2169.7778
$ws.Range("K132").Value = 3073.2309
$ws.Range("L132").Value = 6509.3334
$ws.Range("M132").Value = -543.2309
$ws.Range("N132").Value = -11569.3334

